# Swap the first two comma-separated entries in the "Recorded By" column (G)
# for every data row, leaving any additional entries (e.g. a trailing
# lowercase "system") in place. Rows whose G cell has only a single value
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -notmatch ",") {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $tmp = $parts[0]
    $parts[0] = $parts[1]
    $parts[1] = $tmp

    $cell.Value2 = [string]::Join(", ", $parts)
}
